# Rename the header row to R-friendly (snake_case, lowercase) column names.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A1").Value = "date"
$ws.Range("B1").Value = "distance_in_miles"
$ws.Range("C1").Value = "gasoline_in_gallons"
$ws.Range("D1").Value = "comments"

$ws.Range("A1").Select()
